# Third commit: add a "status" column (Pass/fail) to the login sheet
# and move the active cell selection to B10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Pass"
$ws.Range("C3").Value = "fail"
$ws.Range("C4").Value = "fail"
$ws.Range("C5").Value = "fail"

$ws.Range("B10").Select()
